$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "design": insert a new row for a new question, add supporting links,
# and add extra reference links to a few existing rows.
# ---------------------------------------------------------------------------
$design = $wb.Worksheets.Item("design")

# Remember the existing hyperlink (on the old C12, soon to shift to C13) so
# we can re-attach it after the row insert shuffles everything down.
$existingLinks = @()
foreach ($h in $design.Hyperlinks) {
    $existingLinks += , @($h.Range.Row, $h.Range.Column, $h.Address)
}
foreach ($h in $design.Hyperlinks) {
    $h.Delete()
}

# Insert a brand new row 2 - every row from the old row 2 onward shifts down
# by one.
$design.Rows(2).Insert()

# Re-attach the hyperlink that used to live on C12 - it is now on C13.
foreach ($trip in $existingLinks) {
    $r = $trip[0] + 1
    $c = $trip[1]
    $url = $trip[2]
    $cell = $design.Cells.Item($r, $c)
    $design.Hyperlinks.Add($cell, $url) | Out-Null
}

# New row 2: a brand new interview question.
$design.Range("A2").Value = "设计message queue的读取的API"
$design.Range("B2").Value = 1
$design.Hyperlinks.Add($design.Range("C2"), "https://www.1point3acres.com/bbs/thread-546698-1-1.html") | Out-Null

# Row 5 ("Game of life") gains two reference links.
$design.Hyperlinks.Add($design.Range("C5"), "https://blog.csdn.net/siddontang/article/details/8958323") | Out-Null
$design.Hyperlinks.Add($design.Range("D5"), "https://www.cnblogs.com/wangtao_20/p/3481098.html") | Out-Null

# Row 8 ("log system") is highlighted and gains three reference links.
$design.Range("A8").Font.Color = 255
$design.Range("B8").Font.Color = 255
$design.Range("B8").Value = 2
$design.Hyperlinks.Add($design.Range("C8"), "https://bravenewgeek.com/building-a-distributed-log-from-scratch-part-1-storage-mechanics/") | Out-Null
$design.Hyperlinks.Add($design.Range("D8"), "https://dzone.com/articles/distributed-logging-architecture-for-microservices") | Out-Null
$design.Hyperlinks.Add($design.Range("E8"), "https://www.cnblogs.com/davidwang456/articles/8360274.html") | Out-Null

# Row 15 ("设计Twitter 主页feeds...") is highlighted and gains a reference link.
$design.Range("A15").Font.Color = 255
$design.Range("B15").Value = 2
$design.Hyperlinks.Add($design.Range("C15"), "https://www.1point3acres.com/bbs/thread-545571-1-1.html") | Out-Null

$design.Range("A4").Select()

# ---------------------------------------------------------------------------
# Sheet "freq": bump a couple of counts, add a new data point, and highlight
# a handful of rows (dark-red for the crawler/allocator rows, red for the
# queue/KV-store/setTimeout rows).
# ---------------------------------------------------------------------------
$freq = $wb.Worksheets.Item("freq")

$freq.Range("D4").Value = 10
$freq.Range("A4").Font.Color = 192
$freq.Range("D4").Font.Color = 192
$freq.Range("E4").Font.Color = 192

$freq.Range("B5").Value = 379
$freq.Range("D5").Value = 13
$freq.Range("A5").Font.Color = 192
$freq.Range("B5").Font.Color = 192
$freq.Range("D5").Font.Color = 192

$freq.Range("A13").Font.Color = 255
$freq.Range("D13").Font.Color = 255
$freq.Range("E13").Font.Color = 255
$freq.Range("F13").Font.Color = 255

$freq.Range("A15").Font.Color = 255
$freq.Range("D15").Font.Color = 255
$freq.Range("E15").Font.Color = 255
$freq.Range("F15").Font.Color = 255

$freq.Range("A17").Font.Color = 255
$freq.Range("D17").Font.Color = 255
$freq.Range("E17").Font.Color = 255

$freq.Range("C6").Select()
